# "logboek en verbrokenregels aangepast"
# Add three new logged time entries to the "P6 - Jasper" timesheet and
# switch the active sheet/selection to reflect where the author ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("P6 - Jasper")

# --- New logboek (timesheet) rows -----------------------------------------
# Row 35: Functioneel ontwerp / 26-10-2020 / 30 min / Schermontwerp
$ws.Range("A35").Value = "Functioneel ontwerp"
$ws.Range("B35").Value = 44130
$ws.Range("B34").Copy()
$ws.Range("B35").PasteSpecial(-4122)
$ws.Range("C35").Value = 30
$ws.Range("D35").Value = "Schermontwerp"

# Row 36: Les KBS / 28-10-2020 / 60 min
$ws.Range("A36").Value = "Les KBS"
$ws.Range("B36").Value = 44132
$ws.Range("B34").Copy()
$ws.Range("B36").PasteSpecial(-4122)
$ws.Range("C36").Value = 60

# Row 37: Functioneel ontwerp / 30-10-2020 / 120 min / Domeinmodel
$ws.Range("A37").Value = "Functioneel ontwerp"
$ws.Range("B37").Value = 44134
$ws.Range("B34").Copy()
$ws.Range("B37").PasteSpecial(-4122)
$ws.Range("C37").Value = 120
$ws.Range("D37").Value = "Domeinmodel"

# --- Reflect the author's ending UI state ----------------------------------
# The workbook was left with "P6 - Jasper" as the active tab and the
# selection resting on the next empty logboek row.
$ws.Select() | Out-Null
$ws.Range("A38").Select() | Out-Null

Write-Output "done"
